{"js": "// Replace each three-digit-by-one-digit multiplication prompt in the\n// practice table with its updated problem, per the commit diff.\nconst replacements = {\n  \"698\u00d77=\": \"936\u00d72=\",\n  \"544\u00d75=\": \"683\u00d79=\",\n  \"893\u00d77=\": \"394\u00d76=\",\n  \"585\u00d77=\": \"395\u00d77=\",\n  \"475\u00d77=\": \"856\u00d78=\",\n  \"855\u00d79=\": \"881\u00d75=\",\n  \"331\u00d74=\": \"493\u00d72=\",\n  \"926\u00d73=\": \"234\u00d76=\",\n  \"999\u00d78=\": \"642\u00d72=\",\n  \"364\u00d77=\": \"420\u00d78=\",\n  \"523\u00d75=\": \"973\u00d79=\",\n  \"290\u00d74=\": \"202\u00d73=\",\n  \"598\u00d76=\": \"966\u00d79=\",\n  \"292\u00d76=\": \"599\u00d73=\",\n  \"541\u00d75=\": \"970\u00d75=\",\n  \"111\u00d76=\": \"804\u00d76=\",\n  \"976\u00d77=\": \"951\u00d77=\",\n  \"747\u00d72=\": \"922\u00d74=\",\n  \"616\u00d72=\": \"734\u00d77=\",\n  \"889\u00d76=\": \"308\u00d72=\",\n  \"258\u00d76=\": \"967\u00d76=\",\n  \"478\u00d74=\": \"523\u00d79=\",\n  \"335\u00d77=\": \"230\u00d76=\",\n  \"208\u00d79=\": \"185\u00d75=\",\n  \"265\u00d74=\": \"550\u00d77=\",\n};\n\nfor (const [before, after] of Object.entries(replacements)) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${before}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n", "ps1": "# Update each three-digit-by-one-digit multiplication prompt in the\n# practice table with its updated problem, per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = [ordered]@{\n  '698\u00d77=' = '936\u00d72='\n  '544\u00d75=' = '683\u00d79='\n  '893\u00d77=' = '394\u00d76='\n  '585\u00d77=' = '395\u00d77='\n  '475\u00d77=' = '856\u00d78='\n  '855\u00d79=' = '881\u00d75='\n  '331\u00d74=' = '493\u00d72='\n  '926\u00d73=' = '234\u00d76='\n  '999\u00d78=' = '642\u00d72='\n  '364\u00d77=' = '420\u00d78='\n  '523\u00d75=' = '973\u00d79='\n  '290\u00d74=' = '202\u00d73='\n  '598\u00d76=' = '966\u00d79='\n  '292\u00d76=' = '599\u00d73='\n  '541\u00d75=' = '970\u00d75='\n  '111\u00d76=' = '804\u00d76='\n  '976\u00d77=' = '951\u00d77='\n  '747\u00d72=' = '922\u00d74='\n  '616\u00d72=' = '734\u00d77='\n  '889\u00d76=' = '308\u00d72='\n  '258\u00d76=' = '967\u00d76='\n  '478\u00d74=' = '523\u00d79='\n  '335\u00d77=' = '230\u00d76='\n  '208\u00d79=' = '185\u00d75='\n  '265\u00d74=' = '550\u00d77='\n}\n\nforeach ($before in $replacements.Keys) {\n    $after = $replacements[$before]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $before\n    $find.Replacement.Text = $after\n    $found = $find.Execute(\n        $find.Text, $false, $false, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 2\n    )\n    if (-not $found) {\n        Write-Output \"WARNING: no match for $before\"\n    }\n}\n\n"}
